$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cell C1 changes from "Preventative Health" to "Preventative" for consistency
# with the other single-word motivation headers (Wellness, At Risk, Sick Role, Self Care).
$ws.Range("C1").Value = "Preventative"
